# Insert 3 new rows at position 105, pushing existing rows 105-143 down to 108-146,
# and populate the new rows with the new "Dina" variety records.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(105).Resize(3).Insert()

# Row 105
$ws.Range("A105").Value = 9
$ws.Range("B105").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C105").Value = "Metropolitana"
$ws.Range("D105").Value = 44917
$ws.Range("E105").Value = 13
$ws.Range("F105").Value = "Fruta"
$ws.Range("G105").Value = 100103
$ws.Range("H105").Value = "Frutos de hueso (carozo)"
$ws.Range("I105").Value = 100103003
$ws.Range("J105").Value = "Damasco"
$ws.Range("K105").Value = "Dina"
$ws.Range("L105").Value = "Especial"
$ws.Range("M105").Value = 290
$ws.Range("N105").Value = 19200
$ws.Range("O105").Value = 19200
$ws.Range("P105").Value = 19200
$ws.Range("Q105").Value = "$/caja 16 kilos granel"
$ws.Range("R105").Value = "Región de O'Higgins"
$ws.Range("S105").Value = 1200
$ws.Range("T105").Value = 16

# Row 106
$ws.Range("A106").Value = 9
$ws.Range("B106").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C106").Value = "Metropolitana"
$ws.Range("D106").Value = 44917
$ws.Range("E106").Value = 13
$ws.Range("F106").Value = "Fruta"
$ws.Range("G106").Value = 100103
$ws.Range("H106").Value = "Frutos de hueso (carozo)"
$ws.Range("I106").Value = 100103003
$ws.Range("J106").Value = "Damasco"
$ws.Range("K106").Value = "Dina"
$ws.Range("L106").Value = "Primera"
$ws.Range("M106").Value = 300
$ws.Range("N106").Value = 16000
$ws.Range("O106").Value = 16000
$ws.Range("P106").Value = 16000
$ws.Range("Q106").Value = "$/caja 16 kilos granel"
$ws.Range("R106").Value = "Región de O'Higgins"
$ws.Range("S106").Value = 1000
$ws.Range("T106").Value = 16

# Row 107
$ws.Range("A107").Value = 9
$ws.Range("B107").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C107").Value = "Metropolitana"
$ws.Range("D107").Value = 44917
$ws.Range("E107").Value = 13
$ws.Range("F107").Value = "Fruta"
$ws.Range("G107").Value = 100103
$ws.Range("H107").Value = "Frutos de hueso (carozo)"
$ws.Range("I107").Value = 100103003
$ws.Range("J107").Value = "Damasco"
$ws.Range("K107").Value = "Dina"
$ws.Range("L107").Value = "Segunda"
$ws.Range("M107").Value = 280
$ws.Range("N107").Value = 12800
$ws.Range("O107").Value = 12800
$ws.Range("P107").Value = 12800
$ws.Range("Q107").Value = "$/caja 16 kilos granel"
$ws.Range("R107").Value = "Región de O'Higgins"
$ws.Range("S107").Value = 800
$ws.Range("T107").Value = 16
